# Add a "skill property for each hero" block: insert 3 new columns (Skill1,
# Skill2, Skill3) right before the existing VIPLevel column on the Player
# sheet, mirroring the existing Hero1/Hero2/Hero3 "skill" style columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new columns before column O (VIPLevel). Excel copies formatting
# from the column to the left (the adjacent Skill/Hero block), which matches
# the styles used by the new columns in the target workbook.
$ws.Range("O1:Q1").EntireColumn.Insert()

# New column widths should match the neighbouring Hero columns (K:N).
$ws.Range("O1:Q1").ColumnWidth = $ws.Range("N1").ColumnWidth

# Row 1 (headers)
$ws.Range("O1").Value = "Skill1"
$ws.Range("P1").Value = "Skill2"
$ws.Range("Q1").Value = "Skill3"

# Row 2 (type row)
$ws.Range("O2").Value = "string"
$ws.Range("P2").Value = "string"
$ws.Range("Q2").Value = "string"

# Row 3
$ws.Range("O3").Value = $false
$ws.Range("P3").Value = $false
$ws.Range("Q3").Value = $false

# Row 4
$ws.Range("O4").Value = $true
$ws.Range("P4").Value = $true
$ws.Range("Q4").Value = $true

# Row 5
$ws.Range("O5").Value = $false
$ws.Range("P5").Value = $false
$ws.Range("Q5").Value = $false

# Row 6
$ws.Range("O6").Value = $true
$ws.Range("P6").Value = $true
$ws.Range("Q6").Value = $true

# Row 7
$ws.Range("O7").Value = $false
$ws.Range("P7").Value = $false
$ws.Range("Q7").Value = $false

# Row 8
$ws.Range("O8").Value = $false
$ws.Range("P8").Value = $false
$ws.Range("Q8").Value = $false

# Row 9 (descriptions) - mirrors the Hero1/Hero2/Hero3 labels used by the
# adjacent FightHero block.
$ws.Range("O9").Value = "Hero1"
$ws.Range("P9").Value = "Hero2"
$ws.Range("Q9").Value = "Hero3"

# Restore the frozen-pane view and move the selection the way the authored
# workbook leaves it.
$ws.Range("Q2").Select()
$excel.ActiveWindow.ScrollColumn = 12
